$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates reflecting the refreshed cryptocurrency market data
$updates = @(
    @{ Addr = 'D2'; Value = '28.143.01'; ForceText = $false }
    @{ Addr = 'E2'; Value = '  +0.52%  '; ForceText = $false }
    @{ Addr = 'D3'; Value = '1.876.36'; ForceText = $false }
    @{ Addr = 'E3'; Value = '  +0.85%  '; ForceText = $false }
    @{ Addr = 'D4'; Value = '1.007'; ForceText = $true }
    @{ Addr = 'E4'; Value = '  +0.31%  '; ForceText = $false }
    @{ Addr = 'D5'; Value = '313.58'; ForceText = $true }
    @{ Addr = 'E5'; Value = '  +0.70%  '; ForceText = $false }
    @{ Addr = 'D6'; Value = '1.004'; ForceText = $true }
    @{ Addr = 'E6'; Value = '  +0.16%  '; ForceText = $false }
    @{ Addr = 'D7'; Value = '0.5131'; ForceText = $true }
    @{ Addr = 'E7'; Value = '  +0.26%  '; ForceText = $false }
    @{ Addr = 'D8'; Value = '0.3914'; ForceText = $true }
    @{ Addr = 'E8'; Value = '  +2.75%  '; ForceText = $false }
    @{ Addr = 'D9'; Value = '0.08330'; ForceText = $true }
    @{ Addr = 'E9'; Value = '  +0.52%  '; ForceText = $false }
    @{ Addr = 'D10'; Value = '1.121'; ForceText = $true }
    @{ Addr = 'E10'; Value = '  +1.07%  '; ForceText = $false }
    @{ Addr = 'D11'; Value = '41.46'; ForceText = $true }
    @{ Addr = 'E11'; Value = '  -0.12%  '; ForceText = $false }
    @{ Addr = 'D12'; Value = '6.211'; ForceText = $true }
    @{ Addr = 'E12'; Value = '  +0.11%  '; ForceText = $false }
    @{ Addr = 'D13'; Value = '20.66'; ForceText = $true }
    @{ Addr = 'E13'; Value = '  +0.97%  '; ForceText = $false }
    @{ Addr = 'D14'; Value = '1.881.41'; ForceText = $false }
    @{ Addr = 'E14'; Value = '  +1.15%  '; ForceText = $false }
    @{ Addr = 'D15'; Value = '7.257'; ForceText = $true }
    @{ Addr = 'E15'; Value = '  +0.89%  '; ForceText = $false }
    @{ Addr = 'D16'; Value = '1.005'; ForceText = $true }
    @{ Addr = 'E16'; Value = '  +0.11%  '; ForceText = $false }
    @{ Addr = 'E17'; Value = '  +0.54%  '; ForceText = $false }
    @{ Addr = 'D18'; Value = '91.09'; ForceText = $true }
    @{ Addr = 'E18'; Value = '  +0.61%  '; ForceText = $false }
    @{ Addr = 'D19'; Value = '0.06653'; ForceText = $true }
    @{ Addr = 'E19'; Value = '  +0.63%  '; ForceText = $false }
    @{ Addr = 'D20'; Value = '17.76'; ForceText = $true }
    @{ Addr = 'E20'; Value = '  +0.11%  '; ForceText = $false }
    @{ Addr = 'E21'; Value = '  +0.20%  '; ForceText = $false }
    @{ Addr = 'D22'; Value = '6.030'; ForceText = $true }
    @{ Addr = 'E22'; Value = '  +0.28%  '; ForceText = $false }
    @{ Addr = 'D23'; Value = '28.187.08'; ForceText = $false }
    @{ Addr = 'E23'; Value = '  +0.64%  '; ForceText = $false }
    @{ Addr = 'D24'; Value = '11.14'; ForceText = $true }
    @{ Addr = 'E24'; Value = '  +0.84%  '; ForceText = $false }
    @{ Addr = 'D25'; Value = '2.251'; ForceText = $true }
    @{ Addr = 'E25'; Value = '  +0.86%  '; ForceText = $false }
    @{ Addr = 'B26'; Value = 'WrappedliquidstakedEther2.0'; ForceText = $false }
    @{ Addr = 'C26'; Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'; ForceText = $false }
    @{ Addr = 'D26'; Value = '2.090.55'; ForceText = $false }
    @{ Addr = 'E26'; Value = '  +0.77%  '; ForceText = $false }
    @{ Addr = 'B27'; Value = 'LidoDAOToken'; ForceText = $false }
    @{ Addr = 'C27'; Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; ForceText = $false }
    @{ Addr = 'D27'; Value = '2.493'; ForceText = $true }
    @{ Addr = 'E27'; Value = '  -2.83%  '; ForceText = $false }
    @{ Addr = 'B28'; Value = 'Monero'; ForceText = $false }
    @{ Addr = 'C28'; Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; ForceText = $false }
    @{ Addr = 'D28'; Value = '159.38'; ForceText = $true }
    @{ Addr = 'E28'; Value = '  +1.38%  '; ForceText = $false }
    @{ Addr = 'B29'; Value = 'EthereumClassic'; ForceText = $false }
    @{ Addr = 'C29'; Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; ForceText = $false }
    @{ Addr = 'D29'; Value = '20.61'; ForceText = $true }
    @{ Addr = 'E29'; Value = '  +1.04%  '; ForceText = $false }
    @{ Addr = 'B30'; Value = 'BitcoinCash'; ForceText = $false }
    @{ Addr = 'C30'; Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; ForceText = $false }
    @{ Addr = 'D30'; Value = '125.05'; ForceText = $true }
    @{ Addr = 'E30'; Value = '  +0.19%  '; ForceText = $false }
    @{ Addr = 'B31'; Value = 'Stellar'; ForceText = $false }
    @{ Addr = 'C31'; Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; ForceText = $false }
    @{ Addr = 'D31'; Value = '0.1063'; ForceText = $true }
    @{ Addr = 'E31'; Value = '  +0.09%  '; ForceText = $false }
    @{ Addr = 'B32'; Value = 'ImmutableX'; ForceText = $false }
    @{ Addr = 'C32'; Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; ForceText = $false }
    @{ Addr = 'D32'; Value = '1.038'; ForceText = $true }
    @{ Addr = 'E32'; Value = '  +0.21%  '; ForceText = $false }
    @{ Addr = 'B33'; Value = 'Filecoin'; ForceText = $false }
    @{ Addr = 'C33'; Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; ForceText = $false }
    @{ Addr = 'D33'; Value = '5.843'; ForceText = $true }
    @{ Addr = 'E33'; Value = '  +4.19%  '; ForceText = $false }
    @{ Addr = 'B34'; Value = 'HuobiToken'; ForceText = $false }
    @{ Addr = 'C34'; Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'; ForceText = $false }
    @{ Addr = 'D34'; Value = '3.605'; ForceText = $true }
    @{ Addr = 'E34'; Value = '  -0.03%  '; ForceText = $false }
    @{ Addr = 'B35'; Value = 'FraxShare'; ForceText = $false }
    @{ Addr = 'C35'; Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; ForceText = $false }
    @{ Addr = 'D35'; Value = '9.660'; ForceText = $true }
    @{ Addr = 'E35'; Value = '  +0.82%  '; ForceText = $false }
    @{ Addr = 'B36'; Value = 'VeChain'; ForceText = $false }
    @{ Addr = 'C36'; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; ForceText = $false }
    @{ Addr = 'D36'; Value = '0.02455'; ForceText = $true }
    @{ Addr = 'E36'; Value = '  +1.45%  '; ForceText = $false }
    @{ Addr = 'B37'; Value = 'Hedera'; ForceText = $false }
    @{ Addr = 'C37'; Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'; ForceText = $false }
    @{ Addr = 'D37'; Value = '0.06558'; ForceText = $true }
    @{ Addr = 'E37'; Value = '  +0.17%  '; ForceText = $false }
    @{ Addr = 'B38'; Value = 'Algorand'; ForceText = $false }
    @{ Addr = 'C38'; Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; ForceText = $false }
    @{ Addr = 'D38'; Value = '0.2185'; ForceText = $true }
    @{ Addr = 'E38'; Value = '  +0.37%  '; ForceText = $false }
    @{ Addr = 'B39'; Value = 'ARBITRUM'; ForceText = $false }
    @{ Addr = 'C39'; Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; ForceText = $false }
    @{ Addr = 'D39'; Value = '1.199'; ForceText = $true }
    @{ Addr = 'E39'; Value = '  -0.58%  '; ForceText = $false }
    @{ Addr = 'B40'; Value = 'TheSandbox'; ForceText = $false }
    @{ Addr = 'C40'; Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; ForceText = $false }
    @{ Addr = 'D40'; Value = '0.6495'; ForceText = $true }
    @{ Addr = 'E40'; Value = '  +1.28%  '; ForceText = $false }
    @{ Addr = 'B41'; Value = 'TrustWalletToken'; ForceText = $false }
    @{ Addr = 'C41'; Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; ForceText = $false }
    @{ Addr = 'D41'; Value = '1.229'; ForceText = $true }
    @{ Addr = 'E41'; Value = '  -0.97%  '; ForceText = $false }
    @{ Addr = 'B42'; Value = 'InternetComputer(DFINITY)'; ForceText = $false }
    @{ Addr = 'C42'; Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; ForceText = $false }
    @{ Addr = 'D42'; Value = '4.986'; ForceText = $true }
    @{ Addr = 'E42'; Value = '  +2.18%  '; ForceText = $false }
    @{ Addr = 'B43'; Value = 'Aptos'; ForceText = $false }
    @{ Addr = 'C43'; Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; ForceText = $false }
    @{ Addr = 'D43'; Value = '11.29'; ForceText = $true }
    @{ Addr = 'E43'; Value = '  +0.46%  '; ForceText = $false }
    @{ Addr = 'B44'; Value = 'Decentraland'; ForceText = $false }
    @{ Addr = 'C44'; Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'; ForceText = $false }
    @{ Addr = 'D44'; Value = '0.6141'; ForceText = $true }
    @{ Addr = 'E44'; Value = '  +0.86%  '; ForceText = $false }
    @{ Addr = 'B45'; Value = 'EnergySwap'; ForceText = $false }
    @{ Addr = 'C45'; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; ForceText = $false }
    @{ Addr = 'D45'; Value = '13.09'; ForceText = $true }
    @{ Addr = 'E45'; Value = '  -0.49%  '; ForceText = $false }
    @{ Addr = 'B46'; Value = 'WEMIXTOKEN'; ForceText = $false }
    @{ Addr = 'C46'; Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'; ForceText = $false }
    @{ Addr = 'D46'; Value = '1.283'; ForceText = $true }
    @{ Addr = 'E46'; Value = '  +0.65%  '; ForceText = $false }
    @{ Addr = 'B47'; Value = 'PancakeSwap'; ForceText = $false }
    @{ Addr = 'C47'; Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'; ForceText = $false }
    @{ Addr = 'D47'; Value = '3.673'; ForceText = $true }
    @{ Addr = 'E47'; Value = '  +0.43%  '; ForceText = $false }
    @{ Addr = 'B48'; Value = 'NEARProtocol'; ForceText = $false }
    @{ Addr = 'C48'; Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; ForceText = $false }
    @{ Addr = 'D48'; Value = '2.019'; ForceText = $true }
    @{ Addr = 'E48'; Value = '  +1.95%  '; ForceText = $false }
    @{ Addr = 'B49'; Value = 'EOS'; ForceText = $false }
    @{ Addr = 'C49'; Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'; ForceText = $false }
    @{ Addr = 'D49'; Value = '1.232'; ForceText = $true }
    @{ Addr = 'E49'; Value = '  +2.40%  '; ForceText = $false }
    @{ Addr = 'B50'; Value = 'Quant'; ForceText = $false }
    @{ Addr = 'C50'; Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; ForceText = $false }
    @{ Addr = 'D50'; Value = '120.59'; ForceText = $true }
    @{ Addr = 'E50'; Value = '  -0.09%  '; ForceText = $false }
    @{ Addr = 'B51'; Value = 'Aave'; ForceText = $false }
    @{ Addr = 'C51'; Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'; ForceText = $false }
    @{ Addr = 'D51'; Value = '78.25'; ForceText = $true }
    @{ Addr = 'E51'; Value = '  -1.17%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Addr)
    if ($u.ForceText) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $u.Value
}

Write-Host "Applied $($updates.Count) cell updates"
